$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values stay text (many look numeric, e.g. "1.009"),
# matching the source data which stores them as plain text/inline strings.
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @{
    'D2' = '27.024.44'
    'E2' = '  -0.98%  '
    'D3' = '1.829.65'
    'E3' = '  -0.09%  '
    'D4' = '1.009'
    'E4' = '  -0.17%  '
    'D5' = '311.51'
    'E5' = '  -0.93%  '
    'E6' = '  -0.15%  '
    'D7' = '0.4650'
    'E7' = '  -1.86%  '
    'D8' = '0.3709'
    'D9' = '0.07392'
    'E9' = '  -0.60%  '
    'D10' = '0.8688'
    'E10' = '  -1.80%  '
    'D11' = '20.02'
    'E11' = '  -2.30%  '
    'E12' = '  +7.61%  '
    'D13' = '1.834.35'
    'E13' = '  -1.57%  '
    'D14' = '6.643'
    'E14' = '  +1.41%  '
    'E15' = '  -1.08%  '
    'D16' = '92.08'
    'E16' = '  -1.79%  '
    'E17' = '  +0.15%  '
    'D18' = '0.000009007'
    'E18' = '  +2.50%  '
    'E19' = '  -0.28%  '
    'D20' = '14.69'
    'E20' = '  -0.49%  '
    'D21' = '27.068.02'
    'E21' = '  -2.21%  '
    'D22' = '5.171'
    'E22' = '  -2.15%  '
    'E23' = '  -0.38%  '
    'D24' = '2.065.46'
    'E24' = '  -2.52%  '
    'D25' = '152.74'
    'E25' = '  +0.66%  '
    'E26' = '  -3.15%  '
    'E27' = '  -1.96%  '
    'D28' = '2.100'
    'E28' = '  -1.65%  '
    'D29' = '5.135'
    'E29' = '  -1.79%  '
    'D30' = '115.74'
    'E30' = '  -1.12%  '
    'D31' = '0.08870'
    'E31' = '  -1.24%  '
    'D32' = '2.981'
    'E32' = '  +1.21%  '
    'B33' = 'Filecoin'
    'C33' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D33' = '4.453'
    'E33' = '  -1.85%  '
    'B34' = 'ImmutableX'
    'C34' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D34' = '0.7289'
    'E34' = '  -2.54%  '
    'E35' = '  -3.42%  '
    'D36' = '2.472'
    'E36' = '  +2.24%  '
    'D37' = '1.079'
    'E37' = '  -1.43%  '
    'D38' = '0.01953'
    'E38' = '  -0.10%  '
    'B39' = 'Hedera'
    'C39' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D39' = '0.05250'
    'E39' = '  -1.61%  '
    'B40' = 'FraxShare'
    'C40' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D40' = '7.385'
    'E40' = '  +2.35%  '
    'E41' = '  -0.59%  '
    'D42' = '0.5181'
    'E42' = '  -2.00%  '
    'D43' = '0.1630'
    'E43' = '  -1.55%  '
    'D44' = '0.8575'
    'E44' = '  -15.07%  '
    'D45' = '8.242'
    'E45' = '  -2.86%  '
    'D46' = '0.4847'
    'E46' = '  -1.12%  '
    'B47' = 'EnergySwap'
    'C47' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D47' = '10.25'
    'E47' = '  -2.05%  '
    'B48' = 'PaxDollar'
    'C48' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    'D48' = '1.008'
    'E48' = '  -0.19%  '
    'E49' = '  -2.23%  '
    'E50' = '  -2.31%  '
    'D51' = '0.06250'
    'E51' = '  -0.78%  '
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

Write-Host "Applied $($updates.Count) cell updates"
